# إضافة حدث جديد في Card5 by HOSSAM at 2025-12-08 12:36:15
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card5")

# Format the new row's text-like columns as Text first so values such as
# "5" and date-like strings are not auto-coerced into numbers/dates,
# matching how the rest of the sheet stores this data (plain text).
$newRow = $ws.Range("A16:O16")
$newRow.NumberFormat = "@"

$ws.Cells.Item(16, 1).Value = "5"
$ws.Cells.Item(16, 12).Value = "11/5/2025"
$ws.Cells.Item(16, 13).Value = "قطع سير كويلر مسنن دبل 700"
$ws.Cells.Item(16, 14).Value = "تم تغير سير  دوبل700(محمد نعيم)"
$ws.Cells.Item(16, 15).Value = "فني"
